# The source workbook contains a student CA (continuous assessment) mark
# sheet with per-component scores in columns D:H (Written Assignment,
# Class Test, Lab Record, Presentation, Project Report) for students in
# rows 2:27. This edit clears out all of those component scores, leaving
# only the identifying columns A:C (Student No, Name, Gender), and updates
# the saved cursor/selection state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# Clear the score columns (D:H) for every student row. This is done as two
# separate contiguous blocks (2:16 and 17:27) because that is how the
# workbook was actually edited - rows 2:16 keep their original row "spans"
# bookkeeping while rows 17:27 get it recalculated down to the remaining
# populated columns (A:C).
$ws.Range("D2:H16").ClearContents() | Out-Null
$ws.Range("D17:H27").ClearContents() | Out-Null

# Restore/update the view state: scroll the window so row 10 is at the top
# and select F18 as the active cell, matching the sheet view saved with the
# workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("F18").Select() | Out-Null
